$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Added Timetable Displaying of Course Lists ---
# Clear the stray "Pass" mark on the top-level item 1. (C2) -- it was
# incorrectly marked, so it is cleared back to blank.
$ws.Range("C2").ClearContents()

# Section 3.h "Ability to retrieve and display the timetable data of a
# public course list" (3 points) -- mark its three sub-items as passed
# and record the 3 points earned.
$ws.Range("C60").Value = "x"
$ws.Range("C61").Value = "x"
$ws.Range("C62").Value = "x"
$ws.Range("D59").Value = 3

# Section 4.b "Edit all aspects of an existing course list" (4 points) --
# mark its sub-items as passed and record the 4 points earned.
$ws.Range("C73").Value = "x"
$ws.Range("C74").Value = "x"
$ws.Range("C75").Value = "x"
$ws.Range("C76").Value = "x"
$ws.Range("C77").Value = "x"
$ws.Range("C78").Value = "x"
$ws.Range("C79").Value = "x"
$ws.Range("C80").Value = "x"
$ws.Range("D71").Value = 4

# Section 4.e "Enforce all required attributes when creating or editing a
# course list" (2 points) -- mark its sub-item as passed and record the
# 2 points earned.
$ws.Range("C93").Value = "x"
$ws.Range("D92").Value = 2

# Move the on-screen selection to reflect where the author was working.
$ws.Range("B31").Select()
